$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Câmera Canon - Dólar)
$ws.Range("D2").Value = 5.5158
$ws.Range("E2").Value = 5515.744842
$ws.Range("G2").Value = 7722.043

# Row 3 (Carro Renault - Euro)
$ws.Range("D3").Value = 6.296031973000001
$ws.Range("E3").Value = 28332.1438785
$ws.Range("G3").Value = 56664.288

# Row 4 (Notebook Dell - Dólar)
$ws.Range("D4").Value = 5.5158
$ws.Range("E4").Value = 4964.164842
$ws.Range("G4").Value = 8439.08

# Row 5 (IPhone - Dólar)
$ws.Range("D5").Value = 5.5158
$ws.Range("E5").Value = 4407.124199999999
$ws.Range("G5").Value = 7492.111

# Row 6 (Carro Fiat - Euro)
$ws.Range("D6").Value = 6.296031973000001
$ws.Range("E6").Value = 18888.095919
$ws.Range("G6").Value = 35887.382

# Row 7 (Celular Xiaomi - Dólar)
$ws.Range("D7").Value = 5.5158
$ws.Range("E7").Value = 2650.231584
$ws.Range("G7").Value = 5300.463

# Row 8 (Joia 20g - Ouro)
$ws.Range("D8").Value = 322.39
$ws.Range("E8").Value = 6447.799999999999
$ws.Range("G8").Value = 7414.97
